$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[60.40922745835378, 66.90687120620268]"
$ws.Range("T2").Value = "[47.80622801826759, 52.2118997151283]"
$ws.Range("L3").Value = "[59.40096122643557, 67.99683815699318]"
$ws.Range("T3").Value = "[46.832554744032514, 52.441154125451796]"
